$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Year of Treatment" column (column B); the remaining
# columns (Daily, 4 to 6 days per week, 2 to 3 days per week,
# Once a week or less, Not used in the last 30 days,
# Not known / missing, Total) shift one column to the left.
$ws.Columns("B:B").Delete()

# Append ".jamais.jamais" to every header title except "Country" (A1).
for ($col = 2; $col -le 8; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $current = $cell.Value2
    $cell.Value = $current + ".jamais.jamais"
}
